$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new test case row (row 17) following the existing pattern used by
# rows 2-16, describing "btree_batch_016" (open cases about txn engine).
$ws.Range("A17").Value = "btree_batch_016"
$ws.Range("B17").Value = "y"
$ws.Range("C17").Value = "批量操作语句16执行"
$ws.Range("D17").Value = "batchsql"
$ws.Range("E17").Value = "SingleTable"
$ws.Range("I17").Value = "btree_batch_sql_016"
$ws.Range("K17").Value = "src/test/resources/io.dingodb.test/testdata/btreecases/batchsql/expectedresult/btree_batch_016.csv"
$ws.Range("J17").Value = "select a.NAME as N from b as a order by a.ID"
$ws.Range("G17").Value = "b"
$ws.Range("N17").Value = "csv_equals"

# Widen column I (Batch_sql) so the new, longer value fits.
$ws.Columns("I").ColumnWidth = 20.91

# Match the cursor position recorded by Excel after the edit.
$null = $ws.Range("G24").Select()
